# Add a "Direction" column (F) to the Key messages sheet indicating whether
# a higher value of the indicator is associated with a negative outcome
# (TRUE) or not (FALSE/no clear direction) for people with disabilities.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header in F1
$ws.Range("F1").Value = "Direction"

# Boolean direction flag per indicator row.
# Rows 2, 3, 19 and 26 are intentionally left blank (no direction assigned).
$direction = @{
    4  = $false
    5  = $false
    6  = $false
    7  = $false
    8  = $false
    9  = $false
    10 = $false
    11 = $true
    12 = $false
    13 = $false
    14 = $false
    15 = $false
    16 = $true
    17 = $false
    18 = $false
    20 = $false
    21 = $false
    22 = $false
    23 = $false
    24 = $false
    25 = $false
    27 = $true
    28 = $true
    29 = $true
    30 = $true
}

foreach ($row in $direction.Keys) {
    $ws.Cells.Item($row, 6).Value = $direction[$row]
}

$ws.Range("F30").Select()
